# Apply weekly CompStat crime-data update for 115th Precinct
# (new reporting week 11/28/2022 - 12/4/2022, Volume 29 Number 48)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-NumericCell([string]$addr, [double]$val) {
    $ws.Range($addr).Value = $val
}

function Set-TextCell([string]$addr, [string]$val) {
    $ws.Range($addr).Value = $val
}

function Set-NumericLookingTextCell([string]$addr, [string]$val) {
    # Force a number-looking string to be stored as text (shared string), mirroring
    # cells elsewhere in the sheet that hold literal '0' / '***.*' placeholders as text
    $ws.Range($addr).Value = "'" + $val
}

Set-TextCell "A8" "Volume 29   Number  48"
Set-TextCell "C9" "Report Covering the Week  11/28/2022  Through  12/4/2022"
Set-NumericCell "M14" 0
Set-NumericLookingTextCell "D15" "0"
Set-TextCell "E15" "***.*"
Set-NumericCell "G15" 5
Set-NumericCell "H15" -80
Set-NumericCell "L15" 28.571428571428
Set-NumericCell "M15" 24.137931034482
Set-NumericCell "N15" 9.090909090909
Set-NumericCell "C16" 3
Set-NumericCell "D16" 9
Set-NumericCell "E16" -66.666666666666
Set-NumericCell "F16" 29
Set-NumericCell "G16" 24
Set-NumericCell "H16" 20.833333333333
Set-NumericCell "I16" 291
Set-NumericCell "J16" 231
Set-NumericCell "K16" 25.974025974026
Set-NumericCell "L16" 44.059405940594
Set-NumericCell "M16" -9.907120743034
Set-NumericCell "N16" -76.069078947368
Set-NumericCell "C17" 4
Set-NumericCell "D17" 11
Set-NumericCell "E17" -63.636363636363
Set-NumericCell "F17" 25
Set-NumericCell "G17" 35
Set-NumericCell "H17" -28.571428571428
Set-NumericCell "I17" 386
Set-NumericCell "J17" 364
Set-NumericCell "K17" 6.043956043956
Set-NumericCell "L17" 36.395759717314
Set-NumericCell "M17" 26.557377049180
Set-NumericCell "N17" -3.740648379052
Set-NumericCell "C18" 4
Set-NumericCell "D18" 4
Set-NumericCell "E18" 0
Set-NumericCell "G18" 13
Set-NumericCell "H18" -23.076923076923
Set-NumericCell "I18" 139
Set-NumericCell "J18" 132
Set-NumericCell "K18" 5.303030303030
Set-NumericCell "L18" -13.664596273291
Set-NumericCell "M18" -49.084249084249
Set-NumericCell "N18" -92.441544317563
Set-NumericCell "C19" 20
Set-NumericCell "D19" 101
Set-NumericCell "E19" -80.198019801980
Set-NumericCell "F19" 69
Set-NumericCell "G19" 212
Set-NumericCell "H19" -67.452830188679
Set-NumericCell "I19" 910
Set-NumericCell "J19" 678
Set-NumericCell "K19" 34.218289085545
Set-NumericCell "L19" 85.714285714285
Set-NumericCell "M19" 95.278969957081
Set-NumericCell "N19" -32.542624166048
Set-NumericCell "C20" 6
Set-NumericCell "D20" 1
Set-NumericCell "E20" 500
Set-NumericCell "F20" 27
Set-NumericCell "G20" 12
Set-NumericCell "H20" 125
Set-NumericCell "I20" 295
Set-NumericCell "J20" 191
Set-NumericCell "K20" 54.450261780104
Set-NumericCell "L20" 53.645833333333
Set-NumericCell "M20" 31.111111111111
Set-NumericCell "N20" -85.925572519084
Set-NumericCell "C21" 37
Set-NumericCell "D21" 126
Set-NumericCell "E21" -70.634920634920
Set-NumericCell "F21" 161
Set-NumericCell "G21" 301
Set-NumericCell "H21" -46.511627906976
Set-NumericCell "I21" 2060
Set-NumericCell "J21" 1631
Set-NumericCell "K21" 26.302881667688
Set-NumericCell "L21" 50.694952450621
Set-NumericCell "M21" 26.847290640394
Set-NumericCell "N21" -70.363976406272
Set-NumericCell "D22" 1
Set-NumericCell "E22" -100
Set-NumericLookingTextCell "F22" "0"
Set-NumericCell "H22" -100
Set-NumericCell "J22" 28
Set-NumericCell "K22" 107.142857142857
Set-NumericCell "L22" 241.176470588235
Set-NumericCell "C24" 44
Set-NumericCell "D24" 24
Set-NumericCell "E24" 83.333333333333
Set-NumericCell "F24" 181
Set-NumericCell "G24" 108
Set-NumericCell "H24" 67.592592592592
Set-NumericCell "I24" 1812
Set-NumericCell "J24" 1277
Set-NumericCell "K24" 41.895066562255
Set-NumericCell "L24" 38.638102524866
Set-NumericCell "M24" 71.916508538899
Set-NumericCell "C25" 18
Set-NumericCell "D25" 12
Set-NumericCell "E25" 50
Set-NumericCell "F25" 69
Set-NumericCell "G25" 66
Set-NumericCell "H25" 4.545454545454
Set-NumericCell "I25" 842
Set-NumericCell "J25" 791
Set-NumericCell "K25" 6.447534766118
Set-NumericCell "L25" 19.263456090651
Set-NumericCell "M25" -2.883506343713
Set-NumericCell "D26" 1
Set-NumericCell "G26" 7
Set-NumericCell "H26" -71.428571428571
Set-NumericCell "J26" 53
Set-NumericCell "K26" -9.433962264150
Set-NumericCell "L26" 6.666666666666
Set-NumericCell "C27" 1
Set-NumericCell "D27" 3
Set-NumericCell "E27" -66.666666666666
Set-NumericCell "F27" 11
Set-NumericCell "G27" 7
Set-NumericCell "H27" 57.142857142857
Set-NumericCell "I27" 105
Set-NumericCell "J27" 99
Set-NumericCell "K27" 6.060606060606
Set-NumericCell "L27" 38.157894736842
Set-NumericCell "N28" -85.714285714285
Set-NumericCell "N29" -86.538461538461
